$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column cells keep their text formatting (avoid Excel auto-numeric coercion)
$priceCells = @("D2","D3","D5","D6","D7","D9","D10","D11","D13","D14","D15","D16","D17","D18","D19","D20","D22","D23","D27","D28","D29","D31","D32","D33","D34","D35","D37","D38","D40","D42","D45","D46","D47","D48","D49")
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "43.104.34"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "2.315.91"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "302.57"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").Value = "99.40"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("D7").Value = "0.508"
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.517"
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("D10").Value = "35.98"
$ws.Range("E10").Value = "  +3.67%  "
$ws.Range("D11").Value = "0.0791"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("E12").Value = "  -1.32%  "
$ws.Range("D13").Value = "17.69"
$ws.Range("E13").Value = "  -2.78%  "
$ws.Range("D14").Value = "6.87"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").Value = "2.677.74"
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("D16").Value = "2.297.56"
$ws.Range("E16").Value = "  -1.29%  "
$ws.Range("D17").Value = "0.793"
$ws.Range("E17").Value = "  -3.61%  "
$ws.Range("D18").Value = "43.032.24"
$ws.Range("D19").Value = "13.19"
$ws.Range("E19").Value = "  +4.40%  "
$ws.Range("D20").Value = "6.18"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "68.17"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").Value = "240.48"
$ws.Range("E23").Value = "  +1.14%  "
$ws.Range("E24").Value = "  -3.06%  "
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").Value = "25.11"
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("D28").Value = "169.12"
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("D29").Value = "9.17"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("E30").Value = "  -2.29%  "
$ws.Range("D31").Value = "33.61"
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D32").Value = "4.94"
$ws.Range("E32").Value = "  +4.67%  "
$ws.Range("D33").Value = "5.17"
$ws.Range("E33").Value = "  +2.34%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").Value = "18.35"
$ws.Range("E35").Value = "  +6.56%  "
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("D37").Value = "0.0694"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").Value = "1.82"
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("D40").Value = "2.76"
$ws.Range("E40").Value = "  -2.50%  "
$ws.Range("E41").Value = "  -0.54%  "
$ws.Range("D42").Value = "1.997.92"
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").Value = "17.46"
$ws.Range("E45").Value = "  -1.97%  "
$ws.Range("D46").Value = "2.84"
$ws.Range("E46").Value = "  -1.41%  "
$ws.Range("D47").Value = "76.81"
$ws.Range("E47").Value = "  +8.98%  "
$ws.Range("D48").Value = "54.77"
$ws.Range("E48").Value = "  -2.32%  "
$ws.Range("D49").Value = "2.545.01"
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("E51").Value = "  +3.36%  "
